$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OfferLetter")
$tbl = $ws.ListObjects.Item("OfferLetterList")

# Insert a new column before the current "Attachments" column (column I)
# so data in column I (and onward) shifts right to column J.
$ws.Columns.Item(9).Insert()

# Give the new column (I) its header text.
$ws.Range("I1").Value = "DateOfHiring"

# Match the new column's width to its left neighbor (Status), as Excel does
# when a column is inserted adjacent to an existing formatted column.
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# Grow the table/ListObject to include the new column (now A1:J7).
$tbl.Resize($ws.Range("A1:J7"))

# Resizing the table can lose track of the header text that was already
# present in the last column (previously "Attachments", shifted from I to J)
# -- re-assert it so the table column name stays correct.
$ws.Range("J1").Value = "Attachments"

# Reflect the user's final selection in the saved view state.
$ws.Range("E29").Select()
